# Updates the 'want-to-go' headcount (column F) across all four sheets.
# Mirrors the upstream gh-pages re-scrape commit 456a3b4: small incremental
# bumps to the viewer-interest counters for each event row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 259
$ws.Range("F3").Value = 884
$ws.Range("F4").Value = 558
$ws.Range("F5").Value = 2301
$ws.Range("F6").Value = 1370
$ws.Range("F7").Value = 125
$ws.Range("F8").Value = 816
$ws.Range("F9").Value = 1175
$ws.Range("F10").Value = 1049
$ws.Range("F11").Value = 3061
$ws.Range("F12").Value = 34
$ws.Range("F15").Value = 625
$ws.Range("F16").Value = 535
$ws.Range("F17").Value = 243
$ws.Range("F18").Value = 617
$ws.Range("F19").Value = 1141
$ws.Range("F20").Value = 1141
$ws.Range("F21").Value = 170
$ws.Range("F22").Value = 545
$ws.Range("F23").Value = 198
$ws.Range("F25").Value = 254
$ws.Range("F26").Value = 654
$ws.Range("F27").Value = 608
$ws.Range("F28").Value = 13
$ws.Range("F29").Value = 844
$ws.Range("F30").Value = 93
$ws.Range("F32").Value = 63
$ws.Range("F33").Value = 1059
$ws.Range("F34").Value = 5076
$ws.Range("F35").Value = 535
$ws.Range("F36").Value = 268
$ws.Range("F37").Value = 143
$ws.Range("F39").Value = 8
$ws.Range("F40").Value = 7

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F24").Value = 52
$ws.Range("F25").Value = 391
$ws.Range("F28").Value = 693
$ws.Range("F38").Value = 448
$ws.Range("F40").Value = 15
$ws.Range("F46").Value = 758
$ws.Range("F47").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 646
$ws.Range("F5").Value = 442
$ws.Range("F6").Value = 427

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 646
$ws.Range("F3").Value = 259
$ws.Range("F4").Value = 442
$ws.Range("F6").Value = 884
$ws.Range("F8").Value = 558
$ws.Range("F9").Value = 2301
$ws.Range("F10").Value = 1370
$ws.Range("F11").Value = 125
$ws.Range("F12").Value = 816
$ws.Range("F13").Value = 1175
$ws.Range("F15").Value = 1049
$ws.Range("F16").Value = 3061
$ws.Range("F17").Value = 34
$ws.Range("F21").Value = 427
$ws.Range("F22").Value = 535
$ws.Range("F23").Value = 243
$ws.Range("F24").Value = 617
$ws.Range("F25").Value = 1141
$ws.Range("F26").Value = 1141
$ws.Range("F27").Value = 170
$ws.Range("F29").Value = 545
$ws.Range("F31").Value = 198
$ws.Range("F32").Value = 254
$ws.Range("F33").Value = 52
$ws.Range("F34").Value = 654
$ws.Range("F35").Value = 608
$ws.Range("F36").Value = 391
$ws.Range("F37").Value = 693
$ws.Range("F38").Value = 844
$ws.Range("F39").Value = 93
$ws.Range("F42").Value = 63
$ws.Range("F43").Value = 1059
$ws.Range("F44").Value = 5076
$ws.Range("F46").Value = 535
$ws.Range("F47").Value = 448
$ws.Range("F48").Value = 448
$ws.Range("F49").Value = 268
$ws.Range("F51").Value = 758
$ws.Range("F52").Value = 7

